$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 464
$ws.Range("I2").Value = 1247
$ws.Range("J2").Value = 5242
$ws.Range("K2").Value = 21
$ws.Range("L2").Value = 1410
$ws.Range("M2").Value = 74
$ws.Range("N2").Value = 895
$ws.Range("O2").Value = 5
$ws.Range("P2").Value = 12
$ws.Range("Q2").Value = 11
$ws.Range("R2").Value = 74
$ws.Range("S2").Value = 550
$ws.Range("T2").Value = 937
$ws.Range("U2").Value = 52
$ws.Range("V2").Value = 8073
$ws.Range("X2").Value = 7974
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 133
$ws.Range("AA2").Value = 71
